$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 71.8
$ws.Range("I11").Value = 71.8
$ws.Range("K11").Value = 71.8
$ws.Range("M11").Value = 68.2

$ws.Range("H19").Value = 437.3
$ws.Range("I19").Value = 481.16666
$ws.Range("K19").Value = 481.16666
$ws.Range("M19").Value = -306.16666

$ws.Range("H32").Value = 8022.9165
$ws.Range("I32").Value = 4069
$ws.Range("K32").Value = 4069
$ws.Range("M32").Value = -3743

$ws.Range("H51").Value = 3333.3333
$ws.Range("J51").Value = 4000
$ws.Range("L51").Value = 4000
$ws.Range("N51").Value = -4968

$ws.Range("H98").Value = 895.8
$ws.Range("I98").Value = 869.75
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 869.75
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 628.25
$ws.Range("N98").Value = -3996

$ws.Range("H113").Value = 1626
$ws.Range("I113").Value = 1626
$ws.Range("K113").Value = 1626
$ws.Range("M113").Value = 1628

$ws.Range("H122").Value = 895.8
$ws.Range("I122").Value = 869.75
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2609.25
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -159.25
$ws.Range("N122").Value = -7900

$ws.Range("H135").Value = 1256.5625
$ws.Range("I135").Value = 856.4167
$ws.Range("J135").Value = 2457
$ws.Range("K135").Value = 7707.7503
$ws.Range("L135").Value = 22113
$ws.Range("M135").Value = -5172.7503
$ws.Range("N135").Value = -27183

$ws.Range("H141").Value = 23927.334
$ws.Range("I141").Value = 891
$ws.Range("J141").Value = 70000
$ws.Range("K141").Value = 2673
$ws.Range("L141").Value = 210000
$ws.Range("M141").Value = 2507
$ws.Range("N141").Value = -220360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 4498
$ws.Range("I11").Value = 4992
$ws.Range("J11").Value = 4004
$ws.Range("K11").Value = 4992
$ws.Range("L11").Value = 4004
$ws.Range("M11").Value = -4848
$ws.Range("N11").Value = -4292

$ws.Range("H61").Value = 1387.2858
$ws.Range("I61").Value = 1262.2
$ws.Range("K61").Value = 1262.2
$ws.Range("M61").Value = -1050.2

$ws.Range("H74").Value = 2461.2856
$ws.Range("I74").Value = 2569
$ws.Range("J74").Value = 1815
$ws.Range("K74").Value = 2569
$ws.Range("L74").Value = 1815
$ws.Range("M74").Value = -1695
$ws.Range("N74").Value = -3563

$ws.Range("H77").Value = 2461.2856
$ws.Range("I77").Value = 2569
$ws.Range("J77").Value = 1815
$ws.Range("K77").Value = 12845
$ws.Range("L77").Value = 9075
$ws.Range("M77").Value = -8477
$ws.Range("N77").Value = -17811

$ws.Range("H122").Value = 5999.5
$ws.Range("I122").Value = 5999
$ws.Range("K122").Value = 17997
$ws.Range("M122").Value = -15547

$ws.Range("H136").Value = 1387.2858
$ws.Range("I136").Value = 1262.2
$ws.Range("K136").Value = 3786.6
$ws.Range("M136").Value = -1236.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1047.2
$ws.Range("J80").Value = 1437.5
$ws.Range("L80").Value = 1437.5
$ws.Range("N80").Value = -3433.5

$ws.Range("H83").Value = 1047.2
$ws.Range("J83").Value = 1437.5
$ws.Range("L83").Value = 7187.5
$ws.Range("N83").Value = -17171.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1002.5
$ws.Range("I12").Value = 1002.5
$ws.Range("K12").Value = 1002.5
$ws.Range("M12").Value = -832.5

$ws.Range("H31").Value = 11700.8125
$ws.Range("I31").Value = 12477.667
$ws.Range("J31").Value = 10702
$ws.Range("K31").Value = 12477.667
$ws.Range("L31").Value = 10702
$ws.Range("M31").Value = -12182.667
$ws.Range("N31").Value = -11292

$ws.Range("H34").Value = 11700.8125
$ws.Range("I34").Value = 12477.667
$ws.Range("J34").Value = 10702
$ws.Range("K34").Value = 12477.667
$ws.Range("L34").Value = 10702
$ws.Range("M34").Value = -12275.667
$ws.Range("N34").Value = -11106

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 18.428572
$ws.Range("J12").Value = 31.666666
$ws.Range("L12").Value = 94.99999800000001
$ws.Range("N12").Value = -440.999998

$ws.Range("H109").Value = 861.2857
$ws.Range("J109").Value = 900
$ws.Range("L109").Value = 2700
$ws.Range("N109").Value = -4780

$ws.Range("H113").Value = 1318.5
$ws.Range("I113").Value = 1987.5
$ws.Range("J113").Value = 649.5
$ws.Range("K113").Value = 5962.5
$ws.Range("L113").Value = 1948.5
$ws.Range("M113").Value = -3792.5
$ws.Range("N113").Value = -6288.5

$ws.Range("H122").Value = 5821.4443
$ws.Range("I122").Value = 7825
$ws.Range("J122").Value = 4819.6665
$ws.Range("K122").Value = 70425
$ws.Range("L122").Value = 43376.9985
$ws.Range("M122").Value = -67975
$ws.Range("N122").Value = -48276.9985

$ws.Range("H129").Value = 2066.6667
$ws.Range("I129").Value = 2600
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 7800
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = -2800
$ws.Range("N129").Value = -13000

$ws.Range("H131").Value = 3016.2778
$ws.Range("I131").Value = 1879.5
$ws.Range("J131").Value = 4437.25
$ws.Range("K131").Value = 5638.5
$ws.Range("L131").Value = 13311.75
$ws.Range("M131").Value = -598.5
$ws.Range("N131").Value = -23391.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2079.375
$ws.Range("I97").Value = 1459.5
$ws.Range("K97").Value = 1459.5
$ws.Range("M97").Value = -963.5

$ws.Range("H122").Value = 8997.362999999999
$ws.Range("I122").Value = 5227.3887
$ws.Range("J122").Value = 25962.25
$ws.Range("K122").Value = 15682.1661
$ws.Range("L122").Value = 77886.75
$ws.Range("M122").Value = -13232.1661
$ws.Range("N122").Value = -82786.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2998
$ws.Range("I93").Value = 3497
$ws.Range("K93").Value = 3497
$ws.Range("M93").Value = -2249

$ws.Range("H122").Value = 9999.5
$ws.Range("I122").Value = 9999
$ws.Range("K122").Value = 29997
$ws.Range("M122").Value = -27547

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 24525
$ws.Range("I132").Value = 29125
$ws.Range("K132").Value = 87375
$ws.Range("M132").Value = -84845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10015002
$ws.Range("I14").Value = 10015002
$ws.Range("K14").Value = 10015002
$ws.Range("M14").Value = -10014834

$ws.Range("H94").Value = 27666.334
$ws.Range("J94").Value = 31500
$ws.Range("L94").Value = 31500
$ws.Range("N94").Value = -33302

$ws.Range("H113").Value = 471
$ws.Range("I113").Value = 471
$ws.Range("K113").Value = 1413
$ws.Range("M113").Value = 757

$ws.Range("H132").Value = 8731.944
$ws.Range("I132").Value = 5464
$ws.Range("K132").Value = 16392
$ws.Range("M132").Value = -13862

$ws.Range("H136").Value = 2021.1428
$ws.Range("I136").Value = 2021.1428
$ws.Range("K136").Value = 6063.428400000001
$ws.Range("M136").Value = -3513.428400000001
